$wb = $excel.ActiveWorkbook

# "ip_address_list" sheet: row 1 ("kartoffelnsalat") flag E1 flips
# from boolean FALSE to numeric 1 (connected device can now have its
# IP changed without admin rights).
$wsList = $wb.Worksheets.Item("ip_address_list")
$wsList.Range("E1").Value = 1

# "ip_adress_fav_list" sheet: append the same device as a new favorite
# row (row 3), also flagged with numeric 1.
$wsFav = $wb.Worksheets.Item("ip_adress_fav_list")
$wsFav.Range("A3").Value = "kartoffelnsalat"
$wsFav.Range("B3").Value = "192.168.10.241"
$wsFav.Range("C3").Value = "255.255.255.0"
$wsFav.Range("D3").Value = "kkgg"
$wsFav.Range("E3").Value = 1
